$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13; this shifts the existing rows 13..108 down to 14..109,
# preserving all of their data (matches the diff's observed "row shifted down by one" pattern).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record's data.
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44532
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100112012
$ws.Cells.Item(13, 7).Value = "Espinaca"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 85
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 8000
$ws.Cells.Item(13, 14).Value = "$/docena de atados"
$ws.Cells.Item(13, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(13, 16).Value = 2667
$ws.Cells.Item(13, 17).Value = 3
$ws.Cells.Item(13, 18).Value = "Hortaliza"
